# Updated numbers for 2023.
# The "year" column (A) on both sheets gets bumped by one year for every
# existing data row, and the now-superfluous last row (which represented
# one year further out than needed) is removed entirely.

$wb = $excel.ActiveWorkbook

# ---- Sheet "Mary": rows 2..32 get A += 1, then row 33 is deleted ----
$wsMary = $wb.Worksheets.Item("Mary")
for ($r = 2; $r -le 32; $r++) {
    $cell = $wsMary.Cells.Item($r, 1)
    $cell.Value2 = $cell.Value2 + 1
}
$wsMary.Rows.Item(33).Delete()

# ---- Sheet "John": rows 2..34 get A += 1, then row 35 (the now-blank
# placeholder row carrying only a year in column A) is deleted, which pulls
# the bare trailing "2057" row up from 36 to 35 untouched ----
$wsJohn = $wb.Worksheets.Item("John")
for ($r = 2; $r -le 34; $r++) {
    $cell = $wsJohn.Cells.Item($r, 1)
    $cell.Value2 = $cell.Value2 + 1
}
$wsJohn.Rows.Item(35).Delete()

# ---- Selection / active-sheet bookkeeping ----
# John's selection is parked at D15 (it is not the active tab afterwards).
$wsJohn.Select()
$wsJohn.Range("D15").Select()

# Mary becomes the active tab, selection parked at D25.
$wsMary.Select()
$wsMary.Range("D25").Select()
